# Update the player roster table on Sheet1.
# Column A = Oyuncu Adı (Player), Column B = Pozisyon (Position), Column C = Takım (Team)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @("Oyuncu Adı", "Pozisyon", "Takım"),
    @("Cade Cunningham", "PG,SG", "Detroit Pistons"),
    @("LaMelo Ball", "PG,SG", "Charlotte Hornets"),
    @("Quentin Grimes", "SG,SF", "Philadelphia 76ers"),
    @("Coby White", "PG,SG", "Chicago Bulls"),
    @("Ausar Thompson", "SF,PF", "Detroit Pistons"),
    @("Andrew Wiggins", "SF,PF", "Miami Heat"),
    @("Onyeka Okongwu", "PF,C", "Atlanta Hawks"),
    @("Malik Monk", "PG,SG,SF", "Sacramento Kings"),
    @("Isaiah Hartenstein", "C", "Oklahoma City Thunder"),
    @("Naz Reid", "PF,C", "Minnesota Timberwolves"),
    @("Jusuf Nurkic", "C", "Charlotte Hornets"),
    @("Guerschon Yabusele", "PF,C", "Philadelphia 76ers"),
    @("Damian Lillard", "PG", "Milwaukee Bucks"),
    @("Derrick White", "PG,SG", "Boston Celtics"),
    @("Gary Trent Jr.", "PG,SG,SF", "Milwaukee Bucks"),
    @("Collin Sexton", "PG,SG", "Utah Jazz"),
    @("Anthony Davis", "PF,C", "Dallas Mavericks"),
    @("Robert Williams III", "C", "Portland Trail Blazers")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
